# Add a new "Netherlands" market test-data sheet, modeled on the existing
# "Spain" sheet (same column widths / styles / merged cells), positioned as
# the last tab and made the active sheet.

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Duplicate the Spain sheet and place the copy right after it (i.e. at the end).
$spain.Copy($null, $spain)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Netherlands"

# Market name / ticket reference for the new country.
$newSheet.Range("B2").Value = "Netherlands Market"
$newSheet.Range("B4").Value = "NGC-3144/T2177"

# The copied sheet inherited Spain's explicit row heights for rows 3-5;
# reset them back to the sheet's default height.
$newSheet.Range("A3:A5").EntireRow.AutoFit()

# Spain is no longer the selected/active tab - clear its lingering
# cell-level selection in favor of a "select all" state.
$spain.Select()
$spain.Cells.Select()

# Make the new sheet the active tab, with B4 selected.
$newSheet.Select()
$newSheet.Range("B4").Select()
